{"js": "// Replace the 15 lattice-multiplication problems in the single 5x3 table\n// with a new set of problems, keeping the existing cell/run formatting\n// (sz=32) intact. Each cell's text is 5 lines joined by vertical-tab\n// (\\u000b), which Word's Office.js API treats as a line/\"soft\" break\n// (rendered as <w:br/> between <w:t> runs).\n\nconst VT = \"\\u000b\"; // soft line break inside a paragraph (-> w:br)\n\nfunction cellText(problem, top, left0, left1) {\n  // problem: \"AB x CD\" -> line2 \"  C    D\", line3 \"A|    |\", line4 \"B|    |\"\n  return problem + VT + top + VT + \"  ----\" + VT + left0 + VT + left1;\n}\n\n// New problems laid out row-major, 3 columns x 5 rows (matches tblGrid).\nconst newCells = [\n  cellText(\"86 x 28\", \"  2    8\", \"8|    |\", \"6|    |\"),\n  cellText(\"61 x 50\", \"  5    0\", \"6|    |\", \"1|    |\"),\n  cellText(\"13 x 81\", \"  8    1\", \"1|    |\", \"3|    |\"),\n\n  cellText(\"91 x 23\", \"  2    3\", \"9|    |\", \"1|    |\"),\n  cellText(\"39 x 45\", \"  4    5\", \"3|    |\", \"9|    |\"),\n  cellText(\"33 x 64\", \"  6    4\", \"3|    |\", \"3|    |\"),\n\n  cellText(\"22 x 64\", \"  6    4\", \"2|    |\", \"2|    |\"),\n  cellText(\"28 x 41\", \"  4    1\", \"2|    |\", \"8|    |\"),\n  cellText(\"19 x 48\", \"  4    8\", \"1|    |\", \"9|    |\"),\n\n  cellText(\"82 x 67\", \"  6    7\", \"8|    |\", \"2|    |\"),\n  cellText(\"55 x 89\", \"  8    9\", \"5|    |\", \"5|    |\"),\n  cellText(\"93 x 86\", \"  8    6\", \"9|    |\", \"3|    |\"),\n\n  cellText(\"59 x 96\", \"  9    6\", \"5|    |\", \"9|    |\"),\n  cellText(\"55 x 82\", \"  8    2\", \"5|    |\", \"5|    |\"),\n  cellText(\"25 x 45\", \"  4    5\", \"2|    |\", \"5|    |\"),\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst columnCount = table.values.length > 0 ? table.values[0].length : 0;\n\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    para.insertText(newCells[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 15 lattice-multiplication problems in the single 5x3 table\n# with a new set of problems, keeping the existing cell/run formatting\n# (sz=32) intact. Each cell's text has 5 lines separated by a vertical-tab\n# char (chr(11)), which Word represents as <w:br/> between <w:t> runs.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$VT = [char]11\n\n# New problems laid out row-major, 3 columns x 5 rows (matches tblGrid).\n$newCells = @(\n  \"86 x 28$VT  2    8$VT  ----${VT}8|    |${VT}6|    |\",\n  \"61 x 50$VT  5    0$VT  ----${VT}6|    |${VT}1|    |\",\n  \"13 x 81$VT  8    1$VT  ----${VT}1|    |${VT}3|    |\",\n\n  \"91 x 23$VT  2    3$VT  ----${VT}9|    |${VT}1|    |\",\n  \"39 x 45$VT  4    5$VT  ----${VT}3|    |${VT}9|    |\",\n  \"33 x 64$VT  6    4$VT  ----${VT}3|    |${VT}3|    |\",\n\n  \"22 x 64$VT  6    4$VT  ----${VT}2|    |${VT}2|    |\",\n  \"28 x 41$VT  4    1$VT  ----${VT}2|    |${VT}8|    |\",\n  \"19 x 48$VT  4    8$VT  ----${VT}1|    |${VT}9|    |\",\n\n  \"82 x 67$VT  6    7$VT  ----${VT}8|    |${VT}2|    |\",\n  \"55 x 89$VT  8    9$VT  ----${VT}5|    |${VT}5|    |\",\n  \"93 x 86$VT  8    6$VT  ----${VT}9|    |${VT}3|    |\",\n\n  \"59 x 96$VT  9    6$VT  ----${VT}5|    |${VT}9|    |\",\n  \"55 x 82$VT  8    2$VT  ----${VT}5|    |${VT}5|    |\",\n  \"25 x 45$VT  4    5$VT  ----${VT}2|    |${VT}5|    |\"\n)\n\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $tbl.Cell($r, $c)\n    $rng = $cell.Range\n    # Exclude the trailing cell-end mark (CR + cell marker, 2 chars)\n    # so we only replace the cell's content, not its structure.\n    $rng.End = $rng.End - 2\n    $rng.Text = $newCells[$idx]\n    $idx++\n  }\n}\n"}
